$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column C (shifts old C:G -> F:J)
$ws.Range("C:E").EntireColumn.Insert()

# New header labels for the inserted columns
$ws.Range("C1").Value = "Sentences"
$ws.Range("D1").Value = "Unique_nouns"
$ws.Range("E1").Value = "Unique_verbs"

# New data values for the inserted columns
$ws.Range("C2").Value = 409
$ws.Range("D2").Value = 734
$ws.Range("E2").Value = 782

$ws.Range("C3").Value = 882
$ws.Range("D3").Value = 2068
$ws.Range("E3").Value = 1297

$ws.Range("C4").Value = 31688
$ws.Range("D4").Value = 16698
$ws.Range("E4").Value = 10581
$ws.Range("F4").Value = 4.734268878037985
$ws.Range("G4").Value = 3.825290787898494
$ws.Range("H4").Value = 6.057312252964427
$ws.Range("I4").Value = 5.028447216709195
$ws.Range("J4").Value = 1.794319743472286

$ws.Range("C5").Value = 6236
$ws.Range("D5").Value = 8243
$ws.Range("E5").Value = 7966
